$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 149, shifting existing rows 149-177 down to 150-178
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with the new record's data
$ws.Cells.Item(149, 1).Value = 3
$ws.Cells.Item(149, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(149, 3).Value = "Coquimbo"
$ws.Cells.Item(149, 4).Value = 44711
$ws.Cells.Item(149, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(149, 5).Value = 5
$ws.Cells.Item(149, 6).Value = 100112052
$ws.Cells.Item(149, 7).Value = "Albahaca"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 250
$ws.Cells.Item(149, 11).Value = 4000
$ws.Cells.Item(149, 12).Value = 4500
$ws.Cells.Item(149, 13).Value = 4180
$ws.Cells.Item(149, 14).Value = "$/docena de matas"
$ws.Cells.Item(149, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(149, 16).Value = 697
$ws.Cells.Item(149, 17).Value = 6
$ws.Cells.Item(149, 18).Value = "Hortaliza"
